# Updates the "Estado de Cuenta" (EC) database:
#  - Shifts the list of "Periodo Mora" labels in column E (rows 16-52) forward by
#    one period, so the table now runs from 1806 (oldest) at the top to 2106
#    (newest) at the bottom, instead of 2106 (newest) at the top down to 1806.
#  - The small "Valor Mora" amount (3333) that used to sit on the oldest period
#    (row 16) now belongs to the newest period (row 52), and the previous newest
#    period's amount (70000) now belongs to the oldest period (row 16) -
#    i.e. the two amounts are swapped between rows 16 and 52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    '1806','1807','1808','1809','1810','1811','1812',
    '1901','1902','1903','1904','1905','1906','1907','1908','1909','1910','1911','1912',
    '2001','2002','2003','2004','2005','2006','2007','2008','2009','2010','2011','2012',
    '2101','2102','2103','2104','2105','2106'
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# Swap the "Valor Mora" amounts between the first (row 16) and last (row 52) rows.
$valF16 = $ws.Range("F16").Value2
$valF52 = $ws.Range("F52").Value2
$ws.Range("F16").Value = $valF52
$ws.Range("F52").Value = $valF16
